# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Pil" (C2) becomes "Phil" ---------------------------------
$ws.Range("C2").Value = "Phil"

# --- Row 3: update the FirstName / Company values ----------------------
#   B3: "Golu"   -> "Palak"
#   D3: "UK"     -> "Provar1"
$ws.Range("B3").Value = "Palak"
$ws.Range("D3").Value = "Provar1"

# --- New row 5: repeat the same FirstName / Company values --------------
$ws.Range("B5").Value = "Palak"
$ws.Range("D5").Value = "Provar1"

# --- Restore the cursor position left by the editing session -----------
$ws.Range("D2").Select() | Out-Null
